$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.128.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.05%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.345.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.14%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.07%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.35%  "

# Row 7
$ws.Range("E7").Value = "  -1.58%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.76%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0810"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.74%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.66%  "

# Row 12
$ws.Range("E12").Value = "  +1.05%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.707.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.42%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.68%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.377.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.89%  "

# Row 17
$ws.Range("E17").Value = "  -1.27%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.062.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0902"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.45%  "

# Row 20
$ws.Range("E20").Value = "  -2.57%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.43%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.19%  "

# Row 25
$ws.Range("E25").Value = "  -0.15%  "

# Row 26
$ws.Range("E26").Value = "  -3.63%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.60%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.99%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.68%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.89%  "

# Row 32
$ws.Range("E32").Value = "  +0.00%  "

# Row 33
$ws.Range("E33").Value = "  -2.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0717"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.70%  "

# Row 36
$ws.Range("E36").Value = "  -0.76%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.17%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0985"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.95%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.72%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.86%  "

# Row 41
$ws.Range("E41").Value = "  -0.93%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.968.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.27%  "

# Row 43
$ws.Range("E43").Value = "  -1.58%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0265"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.62%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.78%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.62%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.569.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.06%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "93.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.02%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.32%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.67%  "
